# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: delete C2 entirely, update E2
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 4.566338461218034

# Row 3
$ws.Range("C3").Value = 8.604123301398015
$ws.Range("E3").Value = 8.260999835306748

# Row 4
$ws.Range("E4").Value = 4.862860110364853

# Row 5
$ws.Range("C5").Value = 6.334380382529448

# Row 6
$ws.Range("C6").Value = 3.889938592324382
$ws.Range("E6").Value = 3.430035192100678

# Row 7
$ws.Range("C7").Value = 2.513767348245066

# Row 8
$ws.Range("E8").Value = 2.551173534479356

# Row 10
$ws.Range("C10").Value = 2.321003614014905

# Row 13
$ws.Range("E13").Value = 0.9311475558545279

# Row 14
$ws.Range("E14").Value = 2.152035263856322

# Row 15
$ws.Range("E15").Value = -0.1151106897825049

# Row 18
$ws.Range("C18").Value = -0.1883299148263795

# Row 19
$ws.Range("C19").Value = 3.083905204716264
